# Weekly refresh of "Arveja Verde" price records.
# Existing rows 3-28 get their Fecha / Volumen / Precio mínimo / Precio máximo /
# Precio promedio ponderado / Precio $/Kg values refreshed, and a brand-new
# row 29 (same market/category/quality) is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("D3").Value = 44377
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 28000
$ws.Range("M3").Value = 27000
$ws.Range("P3").Value = 1080

# --- Row 4 ---
$ws.Range("D4").Value = 44350
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 28000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 29000
$ws.Range("P4").Value = 1160

# --- Row 6 ---
$ws.Range("D6").Value = 44349
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 26000
$ws.Range("L6").Value = 28000
$ws.Range("M6").Value = 27000
$ws.Range("P6").Value = 1080

# --- Row 7 ---
$ws.Range("D7").Value = 44364
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("P7").Value = 1160

# --- Row 8 ---
$ws.Range("D8").Value = 44413
$ws.Range("J8").Value = 700
$ws.Range("K8").Value = 26000
$ws.Range("L8").Value = 28000
$ws.Range("M8").Value = 27000
$ws.Range("P8").Value = 1080

# --- Row 9 ---
$ws.Range("D9").Value = 44405
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 26000
$ws.Range("L9").Value = 28000
$ws.Range("M9").Value = 27000
$ws.Range("P9").Value = 1080

# --- Row 10 ---
$ws.Range("D10").Value = 44363
$ws.Range("J10").Value = 240
$ws.Range("K10").Value = 28000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 29000
$ws.Range("P10").Value = 1160

# --- Row 11 ---
$ws.Range("D11").Value = 44343
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 26000
$ws.Range("L11").Value = 28000
$ws.Range("M11").Value = 27000
$ws.Range("P11").Value = 1080

# --- Row 12 ---
$ws.Range("D12").Value = 44406
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 26000
$ws.Range("L12").Value = 28000
$ws.Range("M12").Value = 27000
$ws.Range("P12").Value = 1080

# --- Row 13 ---
$ws.Range("D13").Value = 44385
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 26000
$ws.Range("L13").Value = 28000
$ws.Range("M13").Value = 27000
$ws.Range("P13").Value = 1080

# --- Row 14 ---
$ws.Range("D14").Value = 44371
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 28000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29000
$ws.Range("P14").Value = 1160

# --- Row 15 ---
$ws.Range("D15").Value = 44419
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 27000
$ws.Range("L15").Value = 29000
$ws.Range("M15").Value = 28000
$ws.Range("P15").Value = 1120

# --- Row 16 ---
$ws.Range("D16").Value = 44370
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 27000
$ws.Range("L16").Value = 28000
$ws.Range("M16").Value = 27500
$ws.Range("P16").Value = 1100

# --- Row 17 ---
$ws.Range("D17").Value = 44392
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 26000
$ws.Range("L17").Value = 28000
$ws.Range("M17").Value = 27000
$ws.Range("P17").Value = 1080

# --- Row 18 ---
$ws.Range("D18").Value = 44384
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 26000
$ws.Range("L18").Value = 28000
$ws.Range("M18").Value = 27000
$ws.Range("P18").Value = 1080

# --- Row 19 ---
$ws.Range("D19").Value = 44433
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 28000
$ws.Range("L19").Value = 30000
$ws.Range("M19").Value = 29000
$ws.Range("P19").Value = 1160

# --- Row 20 ---
$ws.Range("D20").Value = 44426
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 28000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 29000
$ws.Range("P20").Value = 1160

# --- Row 21 ---
$ws.Range("D21").Value = 44434
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 28000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 29000
$ws.Range("P21").Value = 1160

# --- Row 22 ---
$ws.Range("D22").Value = 44398
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 26000
$ws.Range("L22").Value = 28000
$ws.Range("M22").Value = 27000
$ws.Range("P22").Value = 1080

# --- Row 23 ---
$ws.Range("D23").Value = 44420
$ws.Range("J23").Value = 700
$ws.Range("K23").Value = 27000
$ws.Range("L23").Value = 29000
$ws.Range("M23").Value = 28000
$ws.Range("P23").Value = 1120

# --- Row 24 ---
$ws.Range("D24").Value = 44427
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 28000
$ws.Range("L24").Value = 30000
$ws.Range("M24").Value = 29000
$ws.Range("P24").Value = 1160

# --- Row 25 ---
$ws.Range("D25").Value = 44441
$ws.Range("J25").Value = 700
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29000
$ws.Range("P25").Value = 1160

# --- Row 26 ---
$ws.Range("D26").Value = 44391
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 26000
$ws.Range("L26").Value = 28000
$ws.Range("M26").Value = 27000
$ws.Range("P26").Value = 1080

# --- Row 27 ---
$ws.Range("D27").Value = 44435
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("P27").Value = 1160

# --- Row 28 ---
$ws.Range("D28").Value = 44412
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 27000
$ws.Range("M28").Value = 26000
$ws.Range("P28").Value = 1040

# --- Row 29 (brand-new record, copied shape of the others) ---
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44399
$ws.Range("D29").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112022
$ws.Range("G29").Value = "Arveja Verde"
$ws.Range("H29").Value = "Perfection"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 26000
$ws.Range("L29").Value = 28000
$ws.Range("M29").Value = 27000
$ws.Range("N29").Value = "`$/malla 25 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 1080
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
